$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 424
$ws.Range("F3").Value = 415
$ws.Range("F4").Value = 2698
$ws.Range("F5").Value = 1318
$ws.Range("F7").Value = 1968
$ws.Range("F8").Value = 567
$ws.Range("F9").Value = 38
$ws.Range("F10").Value = 579
$ws.Range("F11").Value = 271
$ws.Range("F13").Value = 11215
$ws.Range("F14").Value = 6399
$ws.Range("F15").Value = 22
$ws.Range("F21").Value = 883
$ws.Range("F23").Value = 237
$ws.Range("F24").Value = 898
$ws.Range("F25").Value = 3612
$ws.Range("F29").Value = 155
$ws.Range("F30").Value = 302
$ws.Range("F33").Value = 4955
$ws.Range("F35").Value = 1208
$ws.Range("F36").Value = 198
$ws.Range("F37").Value = 377
$ws.Range("F38").Value = 163
$ws.Range("F39").Value = 523

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1
$ws.Range("F9").Value = 137
$ws.Range("F13").Value = 83
$ws.Range("F22").Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8947
$ws.Range("F3").Value = 472
$ws.Range("F4").Value = 1759

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8947
$ws.Range("F3").Value = 472
$ws.Range("F4").Value = 1759
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 415
$ws.Range("F7").Value = 2698
$ws.Range("F10").Value = 1318
$ws.Range("F13").Value = 567
$ws.Range("F14").Value = 38
$ws.Range("F15").Value = 137
$ws.Range("F16").Value = 579
$ws.Range("F17").Value = 271
$ws.Range("F19").Value = 11215
$ws.Range("F21").Value = 6399
$ws.Range("F22").Value = 83
$ws.Range("F23").Value = 22
$ws.Range("F29").Value = 883
$ws.Range("F31").Value = 237
$ws.Range("F32").Value = 898
$ws.Range("F33").Value = 3612
$ws.Range("F35").Value = 155
$ws.Range("F36").Value = 302
$ws.Range("F41").Value = 4955
$ws.Range("F43").Value = 1208
$ws.Range("F44").Value = 199
$ws.Range("F45").Value = 163
$ws.Range("F46").Value = 523
$ws.Range("F48").Value = 37

